$wb = $excel.ActiveWorkbook

# --- Sheet: "All Orders" ---
$orders = $wb.Worksheets.Item("All Orders")
$orders.Range("F2").Value = "Pohe x2, Wheat Chapati x2, Upma x1"
$orders.Range("G2").Value = 120

# --- Sheet: "Daily Summary" ---
$summary = $wb.Worksheets.Item("Daily Summary")
$summary.Range("E2").Value = 405
$summary.Range("G2").Value = 405
